$d = $word.ActiveDocument

# --- Part 1: paragraph 6 - drop proofErr markers, merge " set"/" bit" runs ---
$p6 = $d.Paragraphs(6)
$rng6 = $d.Range($p6.Range.Start, $p6.Range.End)
$xml6 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t xml:space="preserve">The goal of this problem is find the minimum cut for a weighted graph, where the cut is a bisection. The problem uses a binary representation for chromosomes, where </w:t></w:r><w:r><w:t xml:space="preserve">the </w:t></w:r><w:r><w:t xml:space="preserve"> set bit</w:t></w:r><w:r><w:t>s indicate</w:t></w:r><w:r><w:t xml:space="preserve"> one side of the bisection, </w:t></w:r><w:r><w:t>and</w:t></w:r><w:r><w:t xml:space="preserve"> the unset bits indicate t</w:t></w:r><w:r><w:t>he other side</w:t></w:r><w:r><w:t xml:space="preserve">. Chromosomes must have an equal number of 1’s and 0’s in order to </w:t></w:r><w:r><w:t xml:space="preserve">represent a feasible bisection. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng6.InsertXML($xml6)

# --- Part 2: insert new paragraph + extend the fitness-function paragraph ---
$p7 = $d.Paragraphs(7)
$rng7 = $d.Range($p7.Range.Start, $p7.Range.End)
$xml7 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t xml:space="preserve">The graph is read from a file </w:t></w:r><w:r><w:t>that</w:t></w:r><w:r><w:t xml:space="preserve"> lists the vertices </w:t></w:r><w:r><w:t xml:space="preserve">and edges in two separate sets. </w:t></w:r><w:r><w:t>The file is</w:t></w:r><w:r><w:t xml:space="preserve"> parsed and loaded </w:t></w:r><w:r><w:t>in</w:t></w:r><w:r><w:t xml:space="preserve"> a data structure. The structure contains a table of vertices where each vertex is a structure containing the node name and an array of all edges connected to it. </w:t></w:r><w:r><w:t>An</w:t></w:r><w:r><w:t xml:space="preserve"> edge contains pointers to the vertices it connects along with a float</w:t></w:r><w:r><w:t xml:space="preserve">ing point value for its weight. The indices of vertices in the graph table are in the same order </w:t></w:r><w:r><w:t>as the ‘vertices’ (or bits representi</w:t></w:r><w:r><w:t>ng them) are in the chromosome. Any vertex structure can be accessed by a chromosome bit by using its bit position as the graph table’</w:t></w:r><w:r><w:t xml:space="preserve">s array index. </w:t></w:r><w:r><w:t xml:space="preserve">Chromosomes are stored as a dynamic array of 64-bit integers. This allows many operations to process 64-vertices “at once” or in one quad-word for 64-bit processors. If the number of vertices is not a multiple of 64, the last 64-bit integer in the array needs a value to mask out the extra bits. This value is computed at startup. </w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">My fitness function </w:t></w:r><w:r><w:t xml:space="preserve">works by summing the weights of all cut edges. Infeasibles are dealt with by multiplying the difference in number of 1’s and 0’s by 16, and adding the result to the sum of cut edges. </w:t></w:r><w:r><w:t>The number 16 was</w:t></w:r><w:r><w:t xml:space="preserve"> simply chosen through testing.</w:t></w:r><w:r><w:t xml:space="preserve"> Summation requires iterating over all set bits in a chromosome using a find-first-set-bit operator. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng7.InsertXML($xml7)
